$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table Sprint 5")
$ws.Range("P3").Value = 3
$ws.Range("O27").Value = 3
$ws.Range("P27").Value = 1
